$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("command")
$ws.Activate()

# [新規] sub_detectionコマンドを追加 fixed #375
# Insert a new row at 228 (pushes every following row down by one) and
# populate it with the new "sub_detection" command definition.
$ws.Rows(228).Insert()

$ws.Range("A228").Value = "sub_detection"
$ws.Range("B228").Value = "SubDetection"
$ws.Range("C228").Value = "string"
$ws.Range("D228").Value = "double"
$ws.Range("E228").Value = "string"
$ws.Range("G228").Value = "HoI2 DataWiki記載なし"

$ws.Range("A228").Select()
